$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "F33": was "revisione progetti + dire Master DS" -> append "e possibilità tesi su ANN gerarchiche"
$ws.Range("F33").Value = "revisione progetti + dire Master DS e possibilità tesi su ANN gerarchiche"

# "F27": was "hierarchical spotify (partendo dal typing example)" -> reworded
$ws.Range("F27").Value = "hierarchical  (predizione per known and unknown artist, shrinkage, typing exam)"

# "F28": newly filled in (was empty)
$ws.Range("F28").Value = "riprendi hierarchical, regression"

# "F29": newly filled in, copying the "Neutral/Good" highlighted format used in D27
$ws.Range("D27").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Value = "MF: revisione assignment"

# "F30": old note about hierarchical models / ANN thesis removed entirely (folded into F33 above)
$ws.Range("F30").ClearContents()

$excel.CutCopyMode = 0
